$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 1
$ws.Range("B2").Value = 5
$ws.Range("A3").Value = 2
$ws.Range("B3").Value = 6
$ws.Range("A4").Value = 3
$ws.Range("B4").Value = 7
$ws.Range("A5").Value = 4
$ws.Range("B5").Value = 8

$ws.Range("B4").Select()
